$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix individual class-slot swaps (rows 3, 4, 6) before the schedule shift
$ws.Range("F3").Value = "-"
$ws.Range("E4").Value = "MEC-2A-MTRM"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "MCT-2A-MTRM"

# The afternoon schedule (rows 8-14) shifts down by one 50-minute slot,
# starting a new row 9 at 12:20 with "Almoço" that used to be on row 8 (11:30).
# Rebuild rows 8-17 with the corrected time/content grid.
$times = @("11:30", "12:20", "13:00", "13:50", "14:40", "15:30", "15:50", "16:40", "17:30", "18:20")
$rowsData = @(
  @("-", "-", "-", "-", "-"),
  @("Almoço", "Almoço", "Almoço", "Almoço", "Almoço"),
  @("-", "-", "-", "-", "-"),
  @("-", "-", "-", "-", "-"),
  @("-", "-", "-", "-", "-"),
  @("Intervalo", "Intervalo", "Intervalo", "Intervalo", "Intervalo"),
  @("-", "-", "-", "-", "-"),
  @("-", "-", "-", "-", "-"),
  @("-", "-", "-", "-", "-"),
  @("", "", "", "", "")
)

for ($i = 0; $i -lt $times.Length; $i++) {
  $r = 8 + $i
  $ws.Cells.Item($r, 1).Value = $times[$i]
  for ($c = 0; $c -lt 5; $c++) {
    $ws.Cells.Item($r, 2 + $c).Value = $rowsData[$i][$c]
  }
}
